$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.729.09"

# Row 3
$ws.Range("D3").Value = "1.599.43"
$ws.Range("E3").Value = "  +0.20%  "

# Row 5 (D value is a plain number-looking string -> force text with quote prefix)
$ws.Range("D5").Value = "'211.85"
$ws.Range("E5").Value = "  +0.28%  "

# Row 6
$ws.Range("E6").Value = "  +0.13%  "

# Row 7
$ws.Range("E7").Value = "  +0.21%  "

# Row 8
$ws.Range("D8").Value = "'0.0619"
$ws.Range("E8").Value = "  +0.25%  "

# Row 9
$ws.Range("E9").Value = "  +0.23%  "

# Row 10
$ws.Range("D10").Value = "'19.62"
$ws.Range("E10").Value = "  +1.02%  "

# Row 11
$ws.Range("D11").Value = "'0.0846"
$ws.Range("E11").Value = "  +0.54%  "

# Row 12
$ws.Range("D12").Value = "1.823.57"
$ws.Range("E12").Value = "  +0.20%  "

# Row 13
$ws.Range("D13").Value = "1.595.59"
$ws.Range("E13").Value = "  -0.69%  "

# Row 14
$ws.Range("D14").Value = "'4.07"
$ws.Range("E14").Value = "  +1.06%  "

# Row 15
$ws.Range("E15").Value = "  +0.54%  "

# Row 16
$ws.Range("D16").Value = "'65.08"
$ws.Range("E16").Value = "  +0.09%  "

# Row 17
$ws.Range("D17").Value = "0.0₃0740"
$ws.Range("E17").Value = "  -1.42%  "

# Row 19
$ws.Range("D19").Value = "'208.85"
$ws.Range("E19").Value = "  -0.04%  "

# Row 20
$ws.Range("E20").Value = "  +2.51%  "

# Row 21
$ws.Range("E21").Value = "  +0.78%  "

# Row 22
$ws.Range("E22").Value = "  -4.03%  "

# Row 23
$ws.Range("D23").Value = "'9.05"
$ws.Range("E23").Value = "  +1.01%  "

# Row 24
$ws.Range("D24").Value = "'143.78"
$ws.Range("E24").Value = "  +0.67%  "

# Row 25
$ws.Range("E25").Value = "  +0.21%  "

# Row 26
$ws.Range("D26").Value = "'7.13"
$ws.Range("E26").Value = "  +0.27%  "

# Row 27
$ws.Range("E27").Value = "  -0.11%  "

# Row 28
$ws.Range("D28").Value = "'15.38"
$ws.Range("E28").Value = "  +0.40%  "

# Row 29
$ws.Range("E29").Value = "  -1.68%  "

# Row 30
$ws.Range("E30").Value = "  -0.03%  "

# Row 32
$ws.Range("E32").Value = "  +0.65%  "

# Row 33
$ws.Range("D33").Value = "1.281.81"
$ws.Range("E33").Value = "  -0.17%  "

# Row 34
$ws.Range("E34").Value = "  +1.54%  "

# Row 35
$ws.Range("E35").Value = "  +16.86%  "

# Row 36
$ws.Range("E36").Value = "  +0.41%  "

# Row 37
$ws.Range("D37").Value = "'0.591"
$ws.Range("E37").Value = "  -3.88%  "

# Row 38
$ws.Range("E38").Value = "  -1.07%  "

# Row 39
$ws.Range("D39").Value = "'0.825"
$ws.Range("E39").Value = "  +0.02%  "

# Row 40
$ws.Range("E40").Value = "  +0.69%  "

# Row 42
$ws.Range("E42").Value = "  -0.46%  "

# Row 43
$ws.Range("D43").Value = "'62.62"
$ws.Range("E43").Value = "  -0.42%  "

# Row 44
$ws.Range("D44").Value = "1.736.82"
$ws.Range("E44").Value = "  +0.21%  "

# Row 45
$ws.Range("D45").Value = "'90.25"
$ws.Range("E45").Value = "  -0.41%  "

# Row 46
$ws.Range("E46").Value = "  -0.03%  "

# Row 47
$ws.Range("E47").Value = "  +2.17%  "

# Row 48
$ws.Range("E48").Value = "  +0.81%  "

# Row 49
$ws.Range("E49").Value = "  +3.72%  "

# Row 50
$ws.Range("D50").Value = "'1.00"
$ws.Range("E50").Value = "  +0.11%  "

# Row 51
$ws.Range("E51").Value = "  +1.58%  "
